$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

function Replace-Text($old, $new) {
    # Plain literal find & replace across the whole document body.
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: text not found for replace: $old"
    }
    return $ok
}

function Insert-ParaAfter($anchorText, $newText, $style) {
    # Locates the paragraph containing anchorText, inserts a brand-new
    # paragraph immediately after it (inheriting the anchor paragraph's
    # formatting), optionally switches its style, then sets its text.
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "WARNING: anchor not found (after): $anchorText"
        return
    }
    $p = $rng.Paragraphs(1)
    $insertPos = $p.Range.End
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Range($insertPos, $insertPos).Paragraphs(1)
    if ($style) {
        $newPara.Style = $style
    }
    $txtR = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $txtR.Text = $newText
}

function Insert-ParaBefore($anchorText, $newText, $style) {
    # Locates the paragraph containing anchorText, inserts a brand-new
    # paragraph immediately before it (inheriting the anchor paragraph's
    # formatting), optionally switches its style, then sets its text.
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "WARNING: anchor not found (before): $anchorText"
        return
    }
    $p = $rng.Paragraphs(1)
    $anchorStart = $p.Range.Start
    $insR = $d.Range($anchorStart, $anchorStart)
    $insR.InsertParagraphBefore()
    $newPara = $d.Range($anchorStart, $anchorStart).Paragraphs(1)
    if ($style) {
        $newPara.Style = $style
    }
    $txtR = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
    $txtR.Text = $newText
}

# ---------------------------------------------------------------------------
# 1. Add "Author: Majesty Umoye" byline paragraph (centered, gray, 12pt)
#    right after the subtitle paragraph.
# ---------------------------------------------------------------------------

Insert-ParaAfter "ZIP Code 95113 vs Surrounding Areas (2016-2025)" "Author: Majesty Umoye" $null
$byline = $d.Content
$okByline = $byline.Find.Execute("Author: Majesty Umoye", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($okByline) {
    $p = $byline.Paragraphs(1)
    $txtR = $d.Range($p.Range.Start, $p.Range.End - 1)
    $txtR.Font.Color = 6579300
    $txtR.Font.Size = 12
}

# ---------------------------------------------------------------------------
# 2. New "0. Load Masking Correction" subsection at the top of the
#    Executive Summary, before "1. Construction & Commissioning".
# ---------------------------------------------------------------------------

Insert-ParaBefore "1. Construction & Commissioning" "0. Load Masking Correction" "Heading 2"
Insert-ParaAfter "0. Load Masking Correction" "Applied intelligent load masking detection and correction" "List Bullet"
Insert-ParaAfter "Applied intelligent load masking detection and correction" "Identified 23 periods where data center load was hidden in neighboring ZIPs" $null
Insert-ParaAfter "Identified 23 periods where data center load was hidden in neighboring ZIPs" "Correction reveals dramatically higher energy consumption than initially reported" $null
Insert-ParaAfter "Correction reveals dramatically higher energy consumption than initially reported" "True scale of AI/ML workload energy demands now visible" $null

# ---------------------------------------------------------------------------
# 3. Section "1. Construction & Commissioning" bullet edits
# ---------------------------------------------------------------------------

Replace-Text "September 2017: Unprecedented spike to 24.5M kWh" "September 2017: ABSOLUTE PEAK at 24.5M kWh"
Replace-Text "2,200%+ increase suggests initial data center build-out" "Initial data center build-out and infrastructure deployment"
Replace-Text "Followed by intermittent drops during commissioning phase" "Followed by stabilization period as systems came online"

# ---------------------------------------------------------------------------
# 4. Section "2. Steady Operations" bullet edits
# ---------------------------------------------------------------------------

Replace-Text "Relatively stable usage: 8-12M kWh/month" "Relatively stable usage: 10-14M kWh/month"
Replace-Text "Similar to comparison ZIP codes average (10M kWh)" "Standard data center operations"
Replace-Text "Typical data center operations without major expansion" "Pre-AI boom baseline period"
Insert-ParaAfter "Pre-AI boom baseline period" "Typical workloads without intensive ML training" $null

# ---------------------------------------------------------------------------
# 5. Section "3. AI/ML Boom Era" heading + bullet edits
# ---------------------------------------------------------------------------

Replace-Text "3. AI/ML Boom Era (2021-2022)" "3. AI/ML Boom Era (2021-2022) - CRITICAL PERIOD"
Replace-Text "August 2021: Major surge to 17.7M kWh" "August 2021: AI surge begins at 17.7M kWh (GPT-3 adoption)"
Replace-Text "Summer 2022 PEAK: 18.3M kWh (August)" "July-August 2022: SUSTAINED PEAKS at 17.5-18.3M kWh"
Replace-Text "Timeline coincides with GPT-3 adoption, DALL-E 2, and Stable Diffusion" "August 2022: Second-highest peak at 18.3M kWh"
Replace-Text "50-70% increase over 2019 baseline" "Timeline: DALL-E 2 (April 2022), Stable Diffusion (August 2022)"
Insert-ParaAfter "Timeline: DALL-E 2 (April 2022), Stable Diffusion (August 2022)" "Evidence of massive AI model training infrastructure" $null
Insert-ParaAfter "Evidence of massive AI model training infrastructure" "70-80% HIGHER than comparison areas" $null

# ---------------------------------------------------------------------------
# 6. Section "4. Generative AI Explosion" bullet edits
# ---------------------------------------------------------------------------

Replace-Text "March-May 2023: Spikes to 14-15M kWh (GPT-4 launch)" "August 2023: High sustained usage at 17.6M kWh (post-ChatGPT)"
Replace-Text "Sustained high usage throughout 2023-2024" "July 2024: Continued peaks at 17.1M kWh"
Replace-Text "Consistent 13-17M kWh range" "Sustained 16-18M kWh range throughout period"
Replace-Text "Clear correlation with generative AI adoption" "Clear correlation: ChatGPT (Nov 2022), GPT-4 (March 2023)"
Insert-ParaAfter "Clear correlation: ChatGPT (Nov 2022), GPT-4 (March 2023)" "Multi-modal AI models (vision, audio) drive continued high demand" $null
Insert-ParaAfter "Multi-modal AI models (vision, audio) drive continued high demand" "Inference workloads + continued training = sustained energy intensity" $null

# ---------------------------------------------------------------------------
# 7. "Key Findings" -> "Load Masking Correction Applied" paragraph tail edit
#    + new "AI Training Infrastructure Revealed" subsection
# ---------------------------------------------------------------------------

Replace-Text "The excess energy was reallocated back to ZIP 95113, providing a more accurate picture of actual data center consumption." "The excess energy was reallocated back to ZIP 95113, revealing the TRUE scale of AI/ML workload energy consumption that was previously hidden."

Insert-ParaAfter "revealing the TRUE scale of AI/ML workload energy consumption that was previously hidden." "AI Training Infrastructure Revealed" "Heading 2"
Insert-ParaAfter "AI Training Infrastructure Revealed" "Load masking correction exposes the extraordinary energy demands of AI model training. The 2021-2024 period shows sustained peaks of 17-18M kWh per month, with August 2022 reaching 18.3M kWh—the second-highest usage ever recorded. This directly coincides with the training and deployment of GPT-3, DALL-E 2, Stable Diffusion, ChatGPT, and GPT-4. These peaks are 70-80% HIGHER than surrounding comparison areas, demonstrating the massive computational requirements of generative AI." "Normal"

# ---------------------------------------------------------------------------
# 8. "Data Center vs Surrounding Areas" paragraph rewrite
# ---------------------------------------------------------------------------

Replace-Text "The data center ZIP code (95113) shows dramatically different usage patterns compared to surrounding areas, with usage ranging from 90% lower to 180% higher depending on the operational phase. Analysis focused on Commercial, Residential, and Industrial customer classes. Load masking corrections were applied to ensure accurate representation of data center energy consumption." "After load masking correction, the data center ZIP code (95113) shows usage patterns 70-180% higher than surrounding areas during peak AI periods. The corrected data reveals true operational intensity, with consistent monthly usage of 16-18M kWh throughout the AI boom (2021-2024), compared to comparison areas averaging 10M kWh. This represents an unprecedented level of sustained energy consumption driven by AI/ML workloads."

# ---------------------------------------------------------------------------
# 9. "AI Technology Correlation" -> "The AI Timeline: Energy Follows Innovation"
# ---------------------------------------------------------------------------

Replace-Text "AI Technology Correlation" "The AI Timeline: Energy Follows Innovation"
Replace-Text "Strong correlation exists between energy spikes and major AI releases (GPT-3, DALL-E 2, GPT-4, ChatGPT). Training large language models is extremely energy-intensive, requiring exponentially more compute with each generation." "August 2021 (17.7M kWh): GPT-3 widespread adoption begins. July-August 2022 (17.5-18.3M kWh): DALL-E 2 and Stable Diffusion launch—peak AI training period. August 2023 (17.6M kWh): Post-ChatGPT sustained high usage for GPT-4 training. July 2024 (17.1M kWh): Multi-modal AI era continues. The energy data provides a precise timeline of AI development, with each major model release correlating to sustained high energy consumption."

# ---------------------------------------------------------------------------
# 10. "Peak Usage Period (2021-2022)" -> "Peak Usage Period: The AI Boom (2021-2024)"
# ---------------------------------------------------------------------------

Replace-Text "Peak Usage Period (2021-2022)" "Peak Usage Period: The AI Boom (2021-2024)"
Replace-Text "The 2021-2022 period represents peak data center operation, with 50-70% higher energy consumption than surrounding areas. This coincides with the pre-ChatGPT ML boom and widespread AI model training." "The 2021-2024 period represents the AI transformation of data center operations. Unlike traditional steady-state workloads, AI model training creates sustained peaks lasting months. The August 2022 peak of 18.3M kWh coincides with pre-ChatGPT model training, while sustained 2023-2024 usage reflects both continued training and the explosion of inference workloads as ChatGPT and GPT-4 serve millions of users globally."

# ---------------------------------------------------------------------------
# 11. "Lifecycle Phases" -> "Lifecycle Phases With Load Masking Context"
# ---------------------------------------------------------------------------

Replace-Text "Lifecycle Phases" "Lifecycle Phases With Load Masking Context"
Replace-Text "The data clearly shows three distinct phases: (1) Construction/Commissioning (2017), (2) Peak Operations (2018-2024), and (3) Migration/Decommission (2025). Each phase has distinct energy signatures." "The data shows four distinct phases: (1) Construction (Sept 2017: 24.5M kWh absolute peak), (2) Steady Operations (2018-2020: 10-14M kWh), (3) AI Era (2021-2024: sustained 16-18M kWh peaks), and (4) Migration (2025: 700K kWh). Load masking was most prevalent during the AI era, likely due to privacy concerns around revealing the scale of AI infrastructure investment."

# ---------------------------------------------------------------------------
# 12. "Energy Policy Implications" paragraph rewrite
# ---------------------------------------------------------------------------

Replace-Text "Data centers are major energy consumers requiring infrastructure planning. AI/ML workloads create volatile demand patterns that must be accounted for in grid planning and sustainability initiatives." "Load masking correction reveals that AI workloads consume 2-3x more energy than initially apparent. Data centers with AI/ML infrastructure require specialized grid planning, cooling infrastructure, and sustainability strategies. The sustained 17-18M kWh peaks during AI training periods represent energy demands equivalent to powering 15,000+ average homes continuously. Policymakers must account for the exponential growth in AI energy demands when planning grid capacity and renewable energy transitions."

Write-Host "All edits applied."
